$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 2381643.38824391
$ws.Range("E2").Value = 2222955.4972729
$ws.Range("F2").Value = 128535.764409869
$ws.Range("G2").Value = 184.668517112731
$ws.Range("H2").Value = 8.25531196594238

$ws.Range("D3").Value = 2821914.13401452
$ws.Range("E3").Value = 2706256.24919905
$ws.Range("F3").Value = 159403.975752807
$ws.Range("G3").Value = 241.14866900444
$ws.Range("H3").Value = 8.30080509185791

$ws.Range("D4").Value = 3122037.03837648
$ws.Range("E4").Value = 2985050.41841814
$ws.Range("F4").Value = 154130.18182301
$ws.Range("G4").Value = 152.824646234512
$ws.Range("H4").Value = 7.48288106918335

$ws.Range("D5").Value = 3207158.06218067
$ws.Range("E5").Value = 3117791.19862965
$ws.Range("F5").Value = 194508.920347339
$ws.Range("G5").Value = 137.048561811447
$ws.Range("H5").Value = 7.01662373542785

$ws.Range("D6").Value = 3279925.78754841
$ws.Range("E6").Value = 3098092.50487509
$ws.Range("F6").Value = 183063.87272771
$ws.Range("G6").Value = 89.9250009059906
$ws.Range("H6").Value = 6.32664704322814

$ws.Range("D7").Value = 2526654.26392617
$ws.Range("E7").Value = 2447205.65889174
$ws.Range("F7").Value = 214479.017386555
$ws.Range("G7").Value = 69.6399109363555
$ws.Range("H7").Value = 7.23427200317382

$ws.Range("D8").Value = 2829905.23504185
$ws.Range("E8").Value = 2766113.09037919
$ws.Range("F8").Value = 195691.066087126
$ws.Range("G8").Value = 152.655162811279
$ws.Range("H8").Value = 9.60850715637207

$ws.Range("D9").Value = 3011140.61561499
$ws.Range("E9").Value = 3024463.8549595
$ws.Range("F9").Value = 259285.134317752
$ws.Range("G9").Value = 70.4556469917297
$ws.Range("H9").Value = 7.31175780296325

$ws.Range("D10").Value = 3081079.34613573
$ws.Range("E10").Value = 3083956.01693681
$ws.Range("F10").Value = 268673.964345573
$ws.Range("G10").Value = 62.2157700061798
$ws.Range("H10").Value = 7.04303312301635

$ws.Range("D11").Value = 3124470.6658938
$ws.Range("E11").Value = 3077598.91433664
$ws.Range("F11").Value = 263191.501095856
$ws.Range("G11").Value = 62.6546099185943
$ws.Range("H11").Value = 7.0196681022644

# update selection to match target
$ws.Range("A12:XFD18").Select() | Out-Null